$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data reflects reordering each year's block so that Oct, Nov, Dec
# come first, followed by Jan..Sep of that same year.
$data = @(
    @(2, "2014-10", 102.6609, 100.2858, 101.7223),
    @(3, "2014-11", 102.703, 100.2419, 102.4736),
    @(4, "2014-12", 102.6952, 100.247, 111.6992),
    @(5, "2014-01", 101.4589, 100.2257, 101.1956),
    @(6, "2014-02", 101.2153, 100.1673, 101.1492),
    @(7, "2014-03", 100.7677, 100.1676, 101.1492),
    @(8, "2014-04", 100.9247, 100.2232, 101.1492),
    @(9, "2014-05", 100.4312, 100.241, 101.1492),
    @(10, "2014-06", 100.7635, 100.2147, 101.1492),
    @(11, "2014-07", 100.5286, 100.2172, 101.3958),
    @(12, "2014-08", 102.8438, 100.2642, 101.1876),
    @(13, "2014-09", 102.6181, 100.2759, 100.6623),
    @(14, "2015-10", 99.9, 100.1, 111.9),
    @(15, "2015-11", 99.9943, 100.1077, 111.0708),
    @(16, "2015-12", 100.0182, 100.1027, 101.885),
    @(17, "2015-01", 102.5732, 100.1214, 113.6163),
    @(18, "2015-02", 102.5131, 100.1199, 113.6683),
    @(19, "2015-03", 102.1485, 100.1214, 113.6683),
    @(20, "2015-04", 102.3648, 100.0632, 113.6683),
    @(21, "2015-05", 102.1318, 100.1243, 113.6687),
    @(22, "2015-06", 102.0517, 100.1243, 113.6687),
    @(23, "2015-07", 102.1606, 100.1177, 113.3922),
    @(24, "2015-08", 99.856, 100.1098, 113.3922),
    @(25, "2015-09", 99.8696, 100.1098, 113.5798),
    @(26, "2016-10", 100.0, 100.0, 100.0),
    @(27, "2016-11", 99.6, 100.0, 100.0),
    @(28, "2016-12", 99.4, 100.0, 99.9),
    @(29, "2016-01", 98.191, 100.1414, 100.3956),
    @(30, "2016-02", 98.7502, 100.1428, 100.3927),
    @(31, "2016-03", 99.3414, 100.1593, 100.4177),
    @(32, "2016-04", 99.7177, 100.1728, 100.4358),
    @(33, "2016-05", 101.5, 100.1, 100.4),
    @(34, "2016-06", 100.6, 100.1, 100.4),
    @(35, "2016-07", 100.9, 100.1, 100.4),
    @(36, "2016-08", 100.3, 100.0, 100.4),
    @(37, "2016-09", 100.3, 100.0, 100.0),
    @(38, "2017-10", 99.8, 100.0, 99.7),
    @(39, "2017-11", 99.9, 100.0, 99.7),
    @(40, "2017-12", 100.0, 100.0, 99.9),
    @(41, "2017-01", 99.0, 100.0, 99.7),
    @(42, "2017-02", 98.5, 100.0, 100.0),
    @(43, "2017-03", 98.4, 100.0, 100.0),
    @(44, "2017-04", 98.4, 100.0, 100.0),
    @(45, "2017-05", 98.4, 100.0, 100.0),
    @(46, "2017-06", 99.0, 100.0, 100.0),
    @(47, "2017-07", 99.1, 100.0, 100.0),
    @(48, "2017-08", 99.6, 100.0, 100.0),
    @(49, "2017-09", 99.6, 100.0, 100.0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
